# "changes in concise marksheet - Corr/total marks"
# Update the "Marking" row's Right-answer count, and the "Total" row's
# Right-answer count + the Max "correct/total" fraction text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 120
$ws.Range("E12").Value = "120/140"
